$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: EPP Variable Installments T1 scenario cleanup ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Row 2 no longer carries an "Over Due" total in column P
$wsRepay.Range("P2").Clear()

# Rows 3-14 no longer carry a zero in the "Over Due" column O (principal/interest
# totals now land directly in column P)
$wsRepay.Range("O3:O14").Clear()

# --- Transactions sheet: renumber transaction IDs for the new loan scenario ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 1943
$wsTrans.Range("A3").Value = 1942
$wsTrans.Range("A4").Value = 1941
$wsTrans.Range("A5").Value = 1940

# --- View state: restore selections on every sheet, and leave NewLoanInput as
#     the active tab (instead of Transactions) ---
$null = $wb.Worksheets.Item("Summary").Range("D4").Select()
$null = $wsRepay.Range("F11").Select()
$null = $wsTrans.Range("D5").Select()
$null = $wb.Worksheets.Item("NewLoanInput").Range("B2").Select()
